# Automatische test-sync: 2025-07-27 19:50:50
# Appends the new "Testmail #18" row to the Logs sheet and refreshes the
# Dashboard category counts (Bestelling/Levering vs Retour/Terugbetaling)
# to match.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append new row 20 to the Logs sheet -------------------------------
$logs.Range("A20").Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("C20").Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Range("D20").Value = "Bestelling / Levering"
$logs.Range("E20").Value = "Geachte afzender,`nBedankt voor uw e-mail. Het lijkt erop dat deze e-mail bedoeld is als een testmail. Als u daadwerkelijk 200 stuks M8-bouten RVS wenst te bestellen voor Van Dijk, kunt u ons uw bestelgegevens sturen naar het juiste e-mailadres voor bestellingen. Mocht u verdere assistentie nodig hebben of vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F20").Value = "2025-07-27 19:50:12"
$logs.Range("G20").Value = "Ja"
$logs.Range("H20").Value = "Nee"
$logs.Range("I20").Value = "Ja"
$logs.Range("J20").Value = "Nee"

# Re-fit the row height back to the sheet default - the multi-line answer in
# column E would otherwise leave an explicit (custom) row height behind,
# which the source row never had.
$logs.Rows.Item(20).EntireRow.AutoFit()

# --- 2. Extend the conditional-formatting ranges from row 19 to row 20 ----
# (keeps the existing rules / dxf colours - only the applied range grows)
$logs.Range("D2:D19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))
$logs.Range("G2:G19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))
$logs.Range("H2:H19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H20"))
$logs.Range("I2:I19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I20"))
$logs.Range("J2:J19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J20"))

# --- 3. Refresh the Dashboard category-count table -------------------------
# The new row belongs to "Bestelling / Levering", whose count rises from 1
# to 2, tying it with "Retour / Terugbetaling" (also 2). The two rows swap
# places in the ranking, with "Bestelling / Levering" now listed first.
$dashboard.Range("A5").Value = "Bestelling / Levering"
$dashboard.Range("B5").Value = 2
$dashboard.Range("A6").Value = "Retour / Terugbetaling"
$dashboard.Range("B6").Value = 2
